$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price column (D) ---
# Force Text number format on each target cell first so that numeric-looking
# strings (e.g. "1.00", "61.580.29", "0.0000143") are stored as literal text
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.580.29"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.566.15"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.07"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.84"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.42"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.021.48"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.37"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.491.50"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.574.73"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.51"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.79"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.24"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.08"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0833"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.88"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.14"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.39"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.931"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "331.65"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.90"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.43"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.56"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.128.23"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.603"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0545"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.51"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0964"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0238"

# --- Update Volume(1h) column (E) ---
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("E3").Value = "  -5.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -2.66%  "
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("E13").Value = "  -5.71%  "
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("E17").Value = "  -5.50%  "
$ws.Range("E18").Value = "  -5.41%  "
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("E21").Value = "  -5.96%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("E31").Value = "  -5.85%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("E51").Value = "  -3.06%  "
